# Auto-generated Excel COM-interop script to apply the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "lat...4"
$ws.Range("E1").Value = "long...5"
$ws.Range("F1").Value = "lat...6"
$ws.Range("G1").Value = "long...7"
$ws.Range("H1").Value = "lat...8"
$ws.Range("I1").Value = "long...9"
$ws.Range("J1").Value = "lat...10"
$ws.Range("K1").Value = "long...11"
$ws.Range("L1").Value = "lat...12"
$ws.Range("M1").Value = "long...13"
$ws.Range("N1").Value = "lat...14"
$ws.Range("O1").Value = "long...15"
$ws.Range("P1").Value = "lat...16"
$ws.Range("Q1").Value = "long...17"
$ws.Range("F2").Value = 48.86614866335567
$ws.Range("G2").Value = 2.31900834175702
$ws.Range("H2").Value = 48.86614866335567
$ws.Range("I2").Value = 2.31900834175702
$ws.Range("J2").Value = 48.86614866335567
$ws.Range("K2").Value = 2.31900834175702
$ws.Range("L2").Value = 48.87149735714665
$ws.Range("M2").Value = 2.302215602227763
$ws.Range("N2").Value = 48.87149735714665
$ws.Range("O2").Value = 2.302215602227763
$ws.Range("P2").Value = 48.86833055200977
$ws.Range("Q2").Value = 2.312151995261002
$ws.Range("F3").Value = 48.8554814
$ws.Range("G3").Value = 2.3604077
$ws.Range("H3").Value = 48.8554814
$ws.Range("I3").Value = 2.3604077
$ws.Range("J3").Value = 48.8554814
$ws.Range("K3").Value = 2.3604077
$ws.Range("L3").Value = 48.8554814
$ws.Range("M3").Value = 2.3604077
$ws.Range("N3").Value = 48.8554814
$ws.Range("O3").Value = 2.3604077
$ws.Range("P3").Value = 48.8554814
$ws.Range("Q3").Value = 2.3604077
$ws.Range("F4").Value = 45.7588923
$ws.Range("G4").Value = 4.8309221
$ws.Range("H4").Value = 45.7588923
$ws.Range("I4").Value = 4.8309221
$ws.Range("J4").Value = 45.7588923
$ws.Range("K4").Value = 4.8309221
$ws.Range("L4").Value = 45.7588923
$ws.Range("M4").Value = 4.8309221
$ws.Range("N4").Value = 45.7588923
$ws.Range("O4").Value = 4.8309221
$ws.Range("P4").Value = 45.7588923
$ws.Range("Q4").Value = 4.8309221
$ws.Range("F5").Value = 45.7640318
$ws.Range("G5").Value = 4.8356904
$ws.Range("H5").Value = 45.7640318
$ws.Range("I5").Value = 4.8356904
$ws.Range("J5").Value = 45.7640318
$ws.Range("K5").Value = 4.8356904
$ws.Range("L5").Value = 45.7640318
$ws.Range("M5").Value = 4.8356904
$ws.Range("N5").Value = 45.7640318
$ws.Range("O5").Value = 4.8356904
$ws.Range("P5").Value = 45.7640318
$ws.Range("Q5").Value = 4.8356904
